$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Text not found: $find"
    }
    $rng.Text = $replace
}

Replace-Exact 'estudar o job de "pendurar um quadro na parede" ou "criar um escritório atrativo"' `
              'estudar a tarefa (job) de "pendurar um quadro na parede" ou "criar um escritório atrativo"'

Replace-Exact 'você precisaria estudar o job de "criar um buraco". Consequentemente' `
              'você precisaria estudar a tarefa (job) de "criar um buraco". Consequentemente'

Replace-Exact 'estudar o job de "ser bem sucedido na minha profissão" ou "criar um buraco de 60mm"' `
              'estudar a tarefa (job) de "ser bem sucedido na minha profissão" ou "criar um buraco de 60mm"'

Replace-Exact 'você precisaria estudar o job de "pendurar um quadro na parede". Consequentemente' `
              'você precisaria estudar a tarefa (job) de "pendurar um quadro na parede". Consequentemente'

Replace-Exact 'para medir o sucesso ao fazer o job e oportunidades de fazer melhor o job / fazer mais barato.' `
              'para medir o sucesso ao fazer a tarefa (job) e oportunidades de fazê-la melhor/mais barato.'

Replace-Exact 'usar o produto para fazer o job funcional. Mas como' `
              'usar o produto para fazer a tarefa (job) funcional. Mas como'

Replace-Exact 'sobre seu job de comprar (a jornada do comprador)' `
              'sobre sua tarefa de comprar (a jornada do comprador)'

Replace-Exact 'para entender o job que o usuário está tentando fazer (o job funcional central).' `
              'para entender a tarefa que o usuário está tentando fazer (o job funcional central).'

Replace-Exact 'Enquanto estudar o job (jornada) do comprador' `
              'Enquanto estudar a tarefa (jornada) do comprador'
